$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.311.19"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "2.491.02"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'321.16"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").Value = "'108.45"
$ws.Range("E6").Value = "  +3.07%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.534"
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("D10").Value = "'38.74"
$ws.Range("E10").Value = "  +4.49%  "
$ws.Range("D11").Value = "'0.0810"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "'18.34"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").Value = "'7.12"
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("D15").Value = "2.879.18"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "2.484.21"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").Value = "'0.846"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "47.209.43"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "'12.90"
$ws.Range("E19").Value = "  +1.88%  "
$ws.Range("D20").Value = "'6.61"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("E22").Value = "  +12.67%  "
$ws.Range("D23").Value = "'70.28"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("D24").Value = "'245.33"
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("D25").Value = "'2.57"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'25.74"
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("E28").Value = "  +3.61%  "
$ws.Range("D29").Value = "'10.01"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'34.66"
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "'0.137"
$ws.Range("E31").Value = "  +1.92%  "
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").Value = "'20.72"
$ws.Range("E33").Value = "  +5.39%  "
$ws.Range("D34").Value = "'5.34"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("E37").Value = "  +1.53%  "
$ws.Range("D38").Value = "'4.69"
$ws.Range("E38").Value = "  +2.26%  "
$ws.Range("D39").Value = "'2.93"
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("D40").Value = "'23.15"
$ws.Range("E40").Value = "  +8.35%  "
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").Value = "'117.84"
$ws.Range("E43").Value = "  -3.60%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "1.988.84"
$ws.Range("E45").Value = "  +1.23%  "
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("E47").Value = "  -5.87%  "
$ws.Range("D48").Value = "'9.12"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("E49").Value = "  -1.64%  "
$ws.Range("D50").Value = "'5.09"
$ws.Range("D51").Value = "'56.70"
$ws.Range("E51").Value = "  +4.00%  "
